$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2008.091
$ws.Range("I40").Value = 1816.0625
$ws.Range("J40").Value = 2117.8215
$ws.Range("K40").Value = 1816.0625
$ws.Range("L40").Value = 2117.8215
$ws.Range("M40").Value = -1641.0625
$ws.Range("N40").Value = -2467.8215

$ws.Range("H129").Value = 1425404.2
$ws.Range("I129").Value = 260.66666
$ws.Range("J129").Value = 2646956
$ws.Range("K129").Value = 781.9999799999999
$ws.Range("L129").Value = 7940868
$ws.Range("M129").Value = 4218.00002
$ws.Range("N129").Value = -7950868

$ws.Range("H138").Value = 1428.55
$ws.Range("I138").Value = 762.2353000000001
$ws.Range("J138").Value = 1771.803
$ws.Range("K138").Value = 2286.7059
$ws.Range("L138").Value = 5315.409000000001
$ws.Range("M138").Value = 2853.2941
$ws.Range("N138").Value = -15595.409


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 2957
$ws.Range("I96").Value = 2957
$ws.Range("K96").Value = 2957
$ws.Range("M96").Value = -211


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9062.777
$ws.Range("J50").Value = 9062.777
$ws.Range("L50").Value = 9062.777
$ws.Range("N50").Value = -10312.777

$ws.Range("H51").Value = 8844
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 8844
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 8844
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -10316

$ws.Range("H59").Value = 14361.4
$ws.Range("J59").Value = 15701.75
$ws.Range("L59").Value = 15701.75
$ws.Range("N59").Value = -17991.75

$ws.Range("H60").Value = 7700.8335
$ws.Range("I60").Value = 5000
$ws.Range("J60").Value = 8241
$ws.Range("K60").Value = 5000
$ws.Range("L60").Value = 8241
$ws.Range("M60").Value = -4489
$ws.Range("N60").Value = -9263

$ws.Range("H61").Value = 8844
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 8844
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 8844
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -9540

$ws.Range("H68").Value = 15884.875
$ws.Range("J68").Value = 15884.875
$ws.Range("L68").Value = 15884.875
$ws.Range("N68").Value = -17382.875

$ws.Range("H71").Value = 15884.875
$ws.Range("J71").Value = 15884.875
$ws.Range("L71").Value = 47654.625
$ws.Range("N71").Value = -55142.625

$ws.Range("H74").Value = 13526
$ws.Range("J74").Value = 13526
$ws.Range("L74").Value = 13526
$ws.Range("N74").Value = -15274

$ws.Range("H77").Value = 13526
$ws.Range("J77").Value = 13526
$ws.Range("L77").Value = 40578
$ws.Range("N77").Value = -49314

$ws.Range("H132").Value = 1165.7709
$ws.Range("I132").Value = 870.81396
$ws.Range("J132").Value = 3702.4
$ws.Range("K132").Value = 2612.44188
$ws.Range("L132").Value = 11107.2
$ws.Range("M132").Value = -82.44187999999986
$ws.Range("N132").Value = -16167.2


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 333333980
$ws.Range("I75").Value = 913
$ws.Range("J75").Value = 500000500
$ws.Range("K75").Value = 2739
$ws.Range("L75").Value = 1500001500
$ws.Range("M75").Value = -1741
$ws.Range("N75").Value = -1500003496

$ws.Range("H78").Value = 333333980
$ws.Range("I78").Value = 913
$ws.Range("J78").Value = 500000500
$ws.Range("K78").Value = 8217
$ws.Range("L78").Value = 4500004500
$ws.Range("M78").Value = -3225
$ws.Range("N78").Value = -4500014484

$ws.Range("H114").Value = 2254.6667
$ws.Range("I114").Value = 116.5
$ws.Range("J114").Value = 6531
$ws.Range("K114").Value = 349.5
$ws.Range("L114").Value = 19593
$ws.Range("M114").Value = 2904.5
$ws.Range("N114").Value = -26101

$ws.Range("H117").Value = 2668.3845
$ws.Range("I117").Value = 519.3333
$ws.Range("J117").Value = 3313.1
$ws.Range("K117").Value = 1557.9999
$ws.Range("L117").Value = 9939.299999999999
$ws.Range("M117").Value = 1884.0001
$ws.Range("N117").Value = -16823.3

$ws.Range("H131").Value = 8350078.5
$ws.Range("I131").Value = 50092150
$ws.Range("J131").Value = 1663.8363
$ws.Range("K131").Value = 150276450
$ws.Range("L131").Value = 4991.5089
$ws.Range("M131").Value = -150271410
$ws.Range("N131").Value = -15071.5089


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 69337.336
$ws.Range("J20").Value = 69337.336
$ws.Range("L20").Value = 69337.336
$ws.Range("N20").Value = -69827.336

$ws.Range("H22").Value = 304

$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -834

$ws.Range("H70").Value = 4583.6523
$ws.Range("I70").Value = 4172.5454
$ws.Range("K70").Value = 4172.5454
$ws.Range("M70").Value = -3902.5454

$ws.Range("H73").Value = 4583.6523
$ws.Range("I73").Value = 4172.5454
$ws.Range("K73").Value = 4172.5454
$ws.Range("M73").Value = -3236.5454

$ws.Range("H132").Value = 6252698.5
$ws.Range("I132").Value = 8335909.5
$ws.Range("J132").Value = 3065.4
$ws.Range("K132").Value = 25007728.5
$ws.Range("L132").Value = 9196.200000000001
$ws.Range("M132").Value = -25005198.5
$ws.Range("N132").Value = -14256.2


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1419.909
$ws.Range("I46").Value = 1494.875
$ws.Range("K46").Value = 1494.875
$ws.Range("M46").Value = -1306.875

$ws.Range("H93").Value = 1329.6428
$ws.Range("I93").Value = 1270.3846
$ws.Range("J93").Value = 2100
$ws.Range("K93").Value = 1270.3846
$ws.Range("L93").Value = 2100
$ws.Range("M93").Value = -22.38460000000009
$ws.Range("N93").Value = -4596

$ws.Range("H122").Value = 10032.333
$ws.Range("I122").Value = 12687.556
$ws.Range("K122").Value = 38062.66800000001
$ws.Range("M122").Value = -35612.66800000001

$ws.Range("H132").Value = 2164.3901
$ws.Range("I132").Value = 1940.0646
$ws.Range("J132").Value = 2859.8
$ws.Range("K132").Value = 5820.1938
$ws.Range("L132").Value = 8579.400000000001
$ws.Range("M132").Value = -3290.1938
$ws.Range("N132").Value = -13639.4

$ws.Range("H134").Value = 22852.5
$ws.Range("J134").Value = 22852.5
$ws.Range("L134").Value = 22852.5
$ws.Range("N134").Value = -32992.5


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H58").Value = 2875
$ws.Range("I58").Value = 2875
$ws.Range("K58").Value = 2875
$ws.Range("M58").Value = -2567

$ws.Range("H63").Value = 1000
$ws.Range("I63").Value = 1000
$ws.Range("K63").Value = 1000
$ws.Range("M63").Value = -376

$ws.Range("H66").Value = 1000
$ws.Range("I66").Value = 1000
$ws.Range("K66").Value = 3000
$ws.Range("M66").Value = 120

